# Applies the diff: adds two new worksheets (Ejercicio2, Ejercicio3) with
# their data, updates the workbook active-tab selection, and moves the
# "tabSelected" marker from sheet1 to the new last sheet (sheet3).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create "Prueba de escritorio Ejercicio2" right after sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Prueba de escritorio Ejercicio2"

$ws2.Range("B4").Value = "Proceso/SubProcesoLinea(inst)"
$ws2.Range("C4").Value = "PuntosTotales"
$ws2.Range("D4").Value = "T1"
$ws2.Range("E4").Value = "T2"
$ws2.Range("F4").Value = "T3"

$ws2.Range("B5").Value = "1:EJERCICIO21(1)"
$ws2.Range("C5").Value = "<<Variable no inicializada (PUNTOSTOTALES).>>"
$ws2.Range("D5").Value = "<<Variable no inicializada (T1).>>"
$ws2.Range("E5").Value = "<<Variable no inicializada (T2).>>"
$ws2.Range("F5").Value = "<<Variable no inicializada (T3).>>"

$ws2.Range("B6").Value = "1:EJERCICIO22(1)"
$ws2.Range("C6").Value = "<<Variable no inicializada (PUNTOSTOTALES).>>"
$ws2.Range("D6").Value = "<<Variable no inicializada (T1).>>"
$ws2.Range("E6").Value = "<<Variable no inicializada (T2).>>"
$ws2.Range("F6").Value = "<<Variable no inicializada (T3).>>"

$ws2.Range("B7").Value = "1:EJERCICIO25(1)"
$ws2.Range("C7").Value = "<<Variable no inicializada (PUNTOSTOTALES).>>"
$ws2.Range("D7").Value = "<<Variable no inicializada (T1).>>"
$ws2.Range("E7").Value = "<<Variable no inicializada (T2).>>"
$ws2.Range("F7").Value = "<<Variable no inicializada (T3).>>"

$ws2.Range("B8").Value = "1:EJERCICIO26(1)"
$ws2.Range("C8").Value = "<<Variable no inicializada (PUNTOSTOTALES).>>"
$ws2.Range("D8").Value = 50
$ws2.Range("E8").Value = "<<Variable no inicializada (T2).>>"
$ws2.Range("F8").Value = "<<Variable no inicializada (T3).>>"

$ws2.Range("B9").Value = "1:EJERCICIO27(1)"
$ws2.Range("C9").Value = "<<Variable no inicializada (PUNTOSTOTALES).>>"
$ws2.Range("D9").Value = 50
$ws2.Range("E9").Value = 50
$ws2.Range("F9").Value = "<<Variable no inicializada (T3).>>"

$ws2.Range("B10").Value = "1:EJERCICIO29(1)"
$ws2.Range("C10").Value = "<<Variable no inicializada (PUNTOSTOTALES).>>"
$ws2.Range("D10").Value = 50
$ws2.Range("E10").Value = 50
$ws2.Range("F10").Value = 50

$ws2.Range("B11").Value = "1:EJERCICIO210(1)"
$ws2.Range("C11").Value = 150
$ws2.Range("D11").Value = 50
$ws2.Range("E11").Value = 50
$ws2.Range("F11").Value = 50

$ws2.Range("B12").Value = "1:EJERCICIO211(1)"
$ws2.Range("C12").Value = 150
$ws2.Range("D12").Value = 50
$ws2.Range("E12").Value = 50
$ws2.Range("F12").Value = 50

# Column widths matching the "best fit" widths Excel calculated for this
# content (engine only supports 1/6-character granularity when setting
# ColumnWidth, so these are chosen to round-trip as close as possible).
$ws2.Columns.Item(2).ColumnWidth = 28.3072916667
$ws2.Columns.Item(3).ColumnWidth = 43.1666666667
$ws2.Columns.Item(4).ColumnWidth = 30.0221354167
$ws2.Columns.Item(5).ColumnWidth = 30.0221354167
$ws2.Columns.Item(6).ColumnWidth = 30.0221354167

$ws2.Range("G6").Select()

# --- Create "Prueba de escritorio Ejercicio3" right after sheet2 ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Prueba de escritorio Ejercicio3"

$ws3.Range("B4").Value = "Proceso/SubProcesoLinea(inst)"
$ws3.Range("C4").Value = "CantBolsasParaIgualarSaco"
$ws3.Range("D4").Value = "CantidadBolsas"
$ws3.Range("E4").Value = "CostoBolsas"
$ws3.Range("F4").Value = "DiferenciaCosto"
$ws3.Range("G4").Value = "PrecioBolsa"
$ws3.Range("H4").Value = "PrecioSaco"

$ws3.Range("B5").Value = "1:EJERCICIO31(1)"
$ws3.Range("C5").Value = "<<Variable no inicializada (CANTBOLSASPARAIGUALARSACO).>>"
$ws3.Range("D5").Value = "<<Variable no inicializada (CANTIDADBOLSAS).>>"
$ws3.Range("E5").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F5").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G5").Value = "<<Variable no inicializada (PRECIOBOLSA).>>"
$ws3.Range("H5").Value = "<<Variable no inicializada (PRECIOSACO).>>"

$ws3.Range("B6").Value = "1:EJERCICIO33(1)"
$ws3.Range("C6").Value = "<<Variable no inicializada (CANTBOLSASPARAIGUALARSACO).>>"
$ws3.Range("D6").Value = "<<Variable no inicializada (CANTIDADBOLSAS).>>"
$ws3.Range("E6").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F6").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G6").Value = "<<Variable no inicializada (PRECIOBOLSA).>>"
$ws3.Range("H6").Value = "<<Variable no inicializada (PRECIOSACO).>>"

$ws3.Range("B7").Value = "1:EJERCICIO34(1)"
$ws3.Range("C7").Value = "<<Variable no inicializada (CANTBOLSASPARAIGUALARSACO).>>"
$ws3.Range("D7").Value = "<<Variable no inicializada (CANTIDADBOLSAS).>>"
$ws3.Range("E7").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F7").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G7").Value = "<<Variable no inicializada (PRECIOBOLSA).>>"
$ws3.Range("H7").Value = "<<Variable no inicializada (PRECIOSACO).>>"

$ws3.Range("B8").Value = "1:EJERCICIO36(1)"
$ws3.Range("C8").Value = "<<Variable no inicializada (CANTBOLSASPARAIGUALARSACO).>>"
$ws3.Range("D8").Value = "<<Variable no inicializada (CANTIDADBOLSAS).>>"
$ws3.Range("E8").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F8").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G8").Value = "<<Variable no inicializada (PRECIOBOLSA).>>"
$ws3.Range("H8").Value = "<<Variable no inicializada (PRECIOSACO).>>"

$ws3.Range("B9").Value = "1:EJERCICIO37(1)"
$ws3.Range("C9").Value = "<<Variable no inicializada (CANTBOLSASPARAIGUALARSACO).>>"
$ws3.Range("D9").Value = "<<Variable no inicializada (CANTIDADBOLSAS).>>"
$ws3.Range("E9").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F9").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G9").Value = 22
$ws3.Range("H9").Value = "<<Variable no inicializada (PRECIOSACO).>>"

$ws3.Range("B10").Value = "1:EJERCICIO39(1)"
$ws3.Range("C10").Value = "<<Variable no inicializada (CANTBOLSASPARAIGUALARSACO).>>"
$ws3.Range("D10").Value = "<<Variable no inicializada (CANTIDADBOLSAS).>>"
$ws3.Range("E10").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F10").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G10").Value = 22
$ws3.Range("H10").Value = 80

$ws3.Range("B11").Value = "1:EJERCICIO311(1)"
$ws3.Range("C11").Value = 3.6363636364
$ws3.Range("D11").Value = "<<Variable no inicializada (CANTIDADBOLSAS).>>"
$ws3.Range("E11").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F11").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G11").Value = 22
$ws3.Range("H11").Value = 80

$ws3.Range("B12").Value = "1:EJERCICIO313(1)"
$ws3.Range("C12").Value = 3.6363636364
$ws3.Range("D12").Value = "<<Variable no inicializada (CANTIDADBOLSAS).>>"
$ws3.Range("E12").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F12").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G12").Value = 22
$ws3.Range("H12").Value = 80

$ws3.Range("B13").Value = "1:EJERCICIO314(1)"
$ws3.Range("C13").Value = 3.6363636364
$ws3.Range("D13").Value = 12
$ws3.Range("E13").Value = "<<Variable no inicializada (COSTOBOLSAS).>>"
$ws3.Range("F13").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G13").Value = 22
$ws3.Range("H13").Value = 80

$ws3.Range("B14").Value = "1:EJERCICIO316(1)"
$ws3.Range("C14").Value = 3.6363636364
$ws3.Range("D14").Value = 12
$ws3.Range("E14").Value = 264
$ws3.Range("F14").Value = "<<Variable no inicializada (DIFERENCIACOSTO).>>"
$ws3.Range("G14").Value = 22
$ws3.Range("H14").Value = 80

$ws3.Range("B15").Value = "1:EJERCICIO318(1)"
$ws3.Range("C15").Value = 3.6363636364
$ws3.Range("D15").Value = 12
$ws3.Range("E15").Value = 264
$ws3.Range("F15").Value = 184
$ws3.Range("G15").Value = 22
$ws3.Range("H15").Value = 80

$ws3.Range("B16").Value = "1:EJERCICIO320(1)"
$ws3.Range("C16").Value = 3.6363636364
$ws3.Range("D16").Value = 12
$ws3.Range("E16").Value = 264
$ws3.Range("F16").Value = 184
$ws3.Range("G16").Value = 22
$ws3.Range("H16").Value = 80

$ws3.Columns.Item(2).ColumnWidth = 28.3072916667
$ws3.Columns.Item(3).ColumnWidth = 56.8776041667
$ws3.Columns.Item(4).ColumnWidth = 44.0221354167
$ws3.Columns.Item(5).ColumnWidth = 41.0221354167
$ws3.Columns.Item(6).ColumnWidth = 44.7369791667
$ws3.Columns.Item(7).ColumnWidth = 40.3072916667
$ws3.Columns.Item(8).ColumnWidth = 39.4518229167

# sheet3 keeps the "tabSelected" marker (moved from sheet1) plus a wider
# selection rectangle and scrolled viewport.
$ws3.Range("B4:I16").Select()

# sheet1 no longer carries tabSelected - its own selection stays E16.
$ws1.Range("E16").Select()

# Workbook-level: the third tab (index 2, 0-based) is the active one.
$ws3.Activate()
